$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update four individual cells from numeric values to the text "NaN" ---
$ws.Range("CL16").Value2 = "NaN"
$ws.Range("AF127").Value2 = "NaN"
$ws.Range("AF132").Value2 = "NaN"
$ws.Range("AI186").Value2 = "NaN"

# --- Append two new data rows (199 and 200) ---
$row199 = @(44093,758398,2732,102083,66514,250489,27972,5784,4650,7760,8192,17236,3930,22745,30379,7207,9137,14454,13093,17412,14427,3620,2491,9234,26963,13512,10770,56837,1771,841,679,466,547,397,597,2020,4862,37545,8886,2537,44253,1067,22397,1517,9860,1640,1602,7347,1973,955,2491,2661,59451,13631,5484,9377,6179,277,1446,2692,741,2141,9500,9340,10196,14202,1959,893,12673,10255,11932,2482,2000,5167,4445,1816,5499,3275,1941,911,2793,2177,1787,1492,6008,1991,1400,1694,2038,2030,2437,1553,1204,1184,908,3314,1389,924,986,1678,1527,730,855,1246,1514,1405,1467,1171,334,365,783,741,473,543,378,661,746,526,490,373,518,133383,321373,16648,138330,86371,41439,11927)
$row200 = @(44094,765076,2732,103495,66576,252319,28027,5875,4752,7783,8293,17403,3947,22768,30680,7397,9298,14563,13288,17485,14535,3640,2615,9363,27414,13545,10880,57244,1875,938,693,467,631,435,618,2020,4925,37600,8925,2538,44534,1079,22419,1518,9924,1642,1602,7440,1974,955,2491,2661,60616,13642,5589,9427,6273,277,1450,2705,741,2141,9543,9355,10252,14203,1961,893,12793,10378,12186,2548,2028,5250,4501,1885,5513,3296,1961,916,2851,2185,1816,1511,6031,2012,1410,1718,2044,2053,2448,1570,1205,1188,920,3327,1400,925,1009,1681,1533,734,857,1262,1522,1407,1475,1174,334,365,787,743,477,543,380,665,747,526,490,373,518,133966,324061,17003,139798,86969,41998,12054)

$arr199 = New-Object 'object[,]' 1,128
$arr200 = New-Object 'object[,]' 1,128
for ($i = 0; $i -lt 128; $i++) {
    $arr199[0,$i] = $row199[$i]
    $arr200[0,$i] = $row200[$i]
}

$ws.Range("A199:DX199").Value2 = $arr199
$ws.Range("A200:DX200").Value2 = $arr200

# --- Update the active selection to reflect the new last cell ---
$ws.Range("DX199").Select()
